$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.353218793869019
$ws.Range("B1").Value = 4.555501937866211
$ws.Range("C1").Value = 2.495955467224121
$ws.Range("D1").Value = 2.288015365600586
$ws.Range("E1").Value = 1.818939805030823
